# Update imputed values in the RandomForest result sheet.
# These correspond to re-run algorithm output values for columns B and D
# (commit: "Update Name of Algo").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 8.718799999999991
$ws.Range("D4").Value = -7.375299999999998
$ws.Range("B6").Value = 6.266300000000001
$ws.Range("B7").Value = 5.563499999999998
$ws.Range("B8").Value = 6.385599999999999
$ws.Range("D8").Value = -7.820200000000001
$ws.Range("D9").Value = -7.337799999999994
$ws.Range("D12").Value = -6.017399999999996
$ws.Range("B16").Value = 4.9621
$ws.Range("D17").Value = -8.154599999999995
$ws.Range("D18").Value = -9.494199999999987
$ws.Range("D19").Value = -8.666299999999991
$ws.Range("B20").Value = 9.050999999999993
$ws.Range("D20").Value = -8.102299999999991
$ws.Range("B21").Value = 9.240099999999991
$ws.Range("D26").Value = -7.848600000000005
$ws.Range("B28").Value = 6.062400000000008
$ws.Range("B29").Value = 5.201400000000003
$ws.Range("B30").Value = 5.801700000000003
$ws.Range("D31").Value = -7.473899999999994
$ws.Range("B32").Value = 6.764899999999996
$ws.Range("D39").Value = -8.214299999999998
$ws.Range("B40").Value = 9.42499999999999
$ws.Range("D40").Value = -8.698899999999991
$ws.Range("D41").Value = -7.657899999999993
$ws.Range("D42").Value = -8.345899999999993
$ws.Range("D43").Value = -7.468200000000005
$ws.Range("B46").Value = 5.843599999999998
$ws.Range("D47").Value = -7.642199999999998
$ws.Range("D48").Value = -7.558599999999996
$ws.Range("B51").Value = 5.5931
$ws.Range("B52").Value = 5.583699999999997
$ws.Range("D54").Value = -8.027600000000001
$ws.Range("B57").Value = 4.941899999999996
$ws.Range("B59").Value = 4.671600000000001
$ws.Range("B62").Value = 5.105000000000001
$ws.Range("D62").Value = -9.480199999999989
$ws.Range("D63").Value = -6.613199999999997
$ws.Range("D64").Value = -6.930499999999991
$ws.Range("B66").Value = 5.956099999999996
$ws.Range("B73").Value = 8.965599999999997
$ws.Range("B74").Value = 9.004199999999996
$ws.Range("D76").Value = -7.528399999999998
$ws.Range("B77").Value = 8.962900000000005
$ws.Range("D81").Value = -7.645900000000001
$ws.Range("D84").Value = -8.736700000000003
$ws.Range("D89").Value = -8.232699999999999
$ws.Range("B92").Value = 5.036099999999996
$ws.Range("D94").Value = -6.221299999999996
$ws.Range("B100").Value = 5.809600000000002
